# Auto-generated Excel COM-interop script applying the Behemoth_Profits.xlsx edit.
# Updates cached H..N value columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, and clears WVR!M37 which the
# edit removes entirely (no HQ leve for that entry any more).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7109.067
$ws.Range("I62").Value = 3370.913
$ws.Range("J62").Value = 19391.572
$ws.Range("K62").Value = 3370.913
$ws.Range("L62").Value = 19391.572
$ws.Range("M62").Value = -2746.913
$ws.Range("N62").Value = -20639.572
$ws.Range("H65").Value = 7109.067
$ws.Range("I65").Value = 3370.913
$ws.Range("J65").Value = 19391.572
$ws.Range("K65").Value = 16854.565
$ws.Range("L65").Value = 96957.86
$ws.Range("M65").Value = -13734.565
$ws.Range("N65").Value = -103197.86
$ws.Range("H107").Value = 801.93335
$ws.Range("I107").Value = 871.46155
$ws.Range("K107").Value = 871.46155
$ws.Range("M107").Value = 1048.53845
$ws.Range("H116").Value = 4466.1724
$ws.Range("I116").Value = 4153.304
$ws.Range("J116").Value = 5665.5
$ws.Range("K116").Value = 4153.304
$ws.Range("L116").Value = 5665.5
$ws.Range("M116").Value = -711.3040000000001
$ws.Range("N116").Value = -12549.5
$ws.Range("H137").Value = 4911.607
$ws.Range("I137").Value = 2218.3
$ws.Range("J137").Value = 11644.875
$ws.Range("K137").Value = 6654.900000000001
$ws.Range("L137").Value = 34934.625
$ws.Range("M137").Value = -4104.900000000001
$ws.Range("N137").Value = -40034.625
$ws.Range("H138").Value = 2642.6086
$ws.Range("I138").Value = 1824.75
$ws.Range("J138").Value = 2679.7842
$ws.Range("K138").Value = 5474.25
$ws.Range("L138").Value = 8039.3526
$ws.Range("M138").Value = -334.25
$ws.Range("N138").Value = -18319.3526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1196.7693
$ws.Range("I2").Value = 1166.1428
$ws.Range("J2").Value = 1325.4
$ws.Range("K2").Value = 1166.1428
$ws.Range("L2").Value = 1325.4
$ws.Range("M2").Value = -1053.1428
$ws.Range("N2").Value = -1551.4
$ws.Range("H32").Value = 6684131
$ws.Range("I32").Value = 8350650
$ws.Range("J32").Value = 18053.6
$ws.Range("K32").Value = 8350650
$ws.Range("L32").Value = 18053.6
$ws.Range("M32").Value = -8350363
$ws.Range("N32").Value = -18627.6
$ws.Range("H61").Value = 37505860
$ws.Range("I61").Value = 31254374
$ws.Range("K61").Value = 31254374
$ws.Range("M61").Value = -31254162
$ws.Range("H74").Value = 9624433
$ws.Range("I74").Value = 15627304
$ws.Range("K74").Value = 15627304
$ws.Range("M74").Value = -15626430
$ws.Range("H77").Value = 9624433
$ws.Range("I77").Value = 15627304
$ws.Range("K77").Value = 78136520
$ws.Range("M77").Value = -78132152
$ws.Range("H110").Value = 1446.3334
$ws.Range("I110").Value = 1228.8182
$ws.Range("K110").Value = 1228.8182
$ws.Range("M110").Value = 816.1818000000001
$ws.Range("H116").Value = 1196.7693
$ws.Range("I116").Value = 1166.1428
$ws.Range("J116").Value = 1325.4
$ws.Range("K116").Value = 1166.1428
$ws.Range("L116").Value = 1325.4
$ws.Range("M116").Value = 1127.8572
$ws.Range("N116").Value = -5913.4
$ws.Range("H122").Value = 2390.325
$ws.Range("I122").Value = 1601.3793
$ws.Range("K122").Value = 4804.1379
$ws.Range("M122").Value = -2354.1379
$ws.Range("H132").Value = 4418.325
$ws.Range("I132").Value = 2475.8
$ws.Range("K132").Value = 7427.400000000001
$ws.Range("M132").Value = -4897.400000000001
$ws.Range("H136").Value = 37505860
$ws.Range("I136").Value = 31254374
$ws.Range("K136").Value = 93763122
$ws.Range("M136").Value = -93760572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1196.7693
$ws.Range("I3").Value = 1166.1428
$ws.Range("J3").Value = 1325.4
$ws.Range("K3").Value = 1166.1428
$ws.Range("L3").Value = 1325.4
$ws.Range("M3").Value = -1052.1428
$ws.Range("N3").Value = -1553.4
$ws.Range("H98").Value = 67489.5
$ws.Range("J98").Value = 67489.5
$ws.Range("L98").Value = 67489.5
$ws.Range("N98").Value = -73479.5
$ws.Range("H105").Value = 1879.7333
$ws.Range("J105").Value = 2873.4443
$ws.Range("L105").Value = 2873.4443
$ws.Range("N105").Value = -6367.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 659743.3
$ws.Range("I31").Value = 20468.375
$ws.Range("J31").Value = 842393.3
$ws.Range("K31").Value = 20468.375
$ws.Range("L31").Value = 842393.3
$ws.Range("M31").Value = -20173.375
$ws.Range("N31").Value = -842983.3
$ws.Range("H34").Value = 659743.3
$ws.Range("I34").Value = 20468.375
$ws.Range("J34").Value = 842393.3
$ws.Range("K34").Value = 20468.375
$ws.Range("L34").Value = 842393.3
$ws.Range("M34").Value = -20266.375
$ws.Range("N34").Value = -842797.3
$ws.Range("H58").Value = 5402.1665
$ws.Range("I58").Value = 2833.3333
$ws.Range("J58").Value = 7971
$ws.Range("K58").Value = 2833.3333
$ws.Range("L58").Value = 7971
$ws.Range("M58").Value = -2630.3333
$ws.Range("N58").Value = -8377
$ws.Range("H122").Value = 1489.45
$ws.Range("I122").Value = 1374
$ws.Range("K122").Value = 4122
$ws.Range("M122").Value = -1672
$ws.Range("H136").Value = 5402.1665
$ws.Range("I136").Value = 2833.3333
$ws.Range("J136").Value = 7971
$ws.Range("K136").Value = 8499.999899999999
$ws.Range("L136").Value = 23913
$ws.Range("M136").Value = -5949.999899999999
$ws.Range("N136").Value = -29013

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 7143934.5
$ws.Range("J97").Value = 4000
$ws.Range("L97").Value = 12000
$ws.Range("N97").Value = -12992
$ws.Range("H98").Value = 810
$ws.Range("J98").Value = 2000
$ws.Range("L98").Value = 6000
$ws.Range("N98").Value = -8996
$ws.Range("H107").Value = 1323.8182
$ws.Range("J107").Value = 1529.1111
$ws.Range("L107").Value = 4587.3333
$ws.Range("N107").Value = -8427.3333
$ws.Range("H131").Value = 16322
$ws.Range("J131").Value = 18364.223
$ws.Range("L131").Value = 55092.66900000001
$ws.Range("N131").Value = -65172.66900000001
$ws.Range("H133").Value = 2905.8
$ws.Range("I133").Value = 2905.8
$ws.Range("K133").Value = 8717.400000000001
$ws.Range("M133").Value = -3657.400000000001
$ws.Range("H137").Value = 7555.5713
$ws.Range("J137").Value = 9672.25
$ws.Range("L137").Value = 29016.75
$ws.Range("N137").Value = -39216.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 199.16667
$ws.Range("I107").Value = 73.75
$ws.Range("K107").Value = 73.75
$ws.Range("M107").Value = 1846.25
$ws.Range("H108").Value = 130449.25
$ws.Range("J108").Value = 130449.25
$ws.Range("L108").Value = 130449.25
$ws.Range("N108").Value = -138129.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7057.125
$ws.Range("I122").Value = 7009.8
$ws.Range("K122").Value = 21029.4
$ws.Range("M122").Value = -18579.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
